# Update "想去人数" (F column) figures across the four sheets of the
# workbook to match the freshly generated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value2  = 6314
$ws.Range("F3").Value2  = 80
$ws.Range("F5").Value2  = 419
$ws.Range("F6").Value2  = 1401
$ws.Range("F7").Value2  = 3159
$ws.Range("F9").Value2  = 483
$ws.Range("F10").Value2 = 2030
$ws.Range("F11").Value2 = 131
$ws.Range("F13").Value2 = 211
$ws.Range("F14").Value2 = 99
$ws.Range("F15").Value2 = 204
$ws.Range("F16").Value2 = 1018
$ws.Range("F17").Value2 = 383
$ws.Range("F19").Value2 = 130
$ws.Range("F20").Value2 = 3841
$ws.Range("F22").Value2 = 3042
$ws.Range("F24").Value2 = 45
$ws.Range("F25").Value2 = 2582
$ws.Range("F26").Value2 = 4404
$ws.Range("F28").Value2 = 942
$ws.Range("F29").Value2 = 493
$ws.Range("F30").Value2 = 2950
$ws.Range("F31").Value2 = 192
$ws.Range("F32").Value2 = 22
$ws.Range("F33").Value2 = 63
$ws.Range("F34").Value2 = 51
$ws.Range("F35").Value2 = 41
$ws.Range("F36").Value2 = 1064
$ws.Range("F37").Value2 = 1310
$ws.Range("F39").Value2 = 1160
$ws.Range("F40").Value2 = 750
$ws.Range("F41").Value2 = 2
$ws.Range("F42").Value2 = 679
$ws.Range("F43").Value2 = 453
$ws.Range("F44").Value2 = 32
$ws.Range("F45").Value2 = 153
$ws.Range("F46").Value2 = 10
$ws.Range("F47").Value2 = 27
$ws.Range("F48").Value2 = 333
$ws.Range("F49").Value2 = 3637

# --- Sheet "演出" (performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value2 = 924
$ws.Range("F25").Value2 = 23

# --- Sheet "本地生活" (local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value2 = 544

# --- Sheet "全部类型" (all types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value2  = 6315
$ws.Range("F3").Value2  = 80
$ws.Range("F6").Value2  = 419
$ws.Range("F7").Value2  = 1401
$ws.Range("F8").Value2  = 3159
$ws.Range("F9").Value2  = 483
$ws.Range("F11").Value2 = 2030
$ws.Range("F12").Value2 = 131
$ws.Range("F15").Value2 = 211
$ws.Range("F16").Value2 = 924
$ws.Range("F18").Value2 = 99
$ws.Range("F19").Value2 = 204
$ws.Range("F20").Value2 = 1018
$ws.Range("F22").Value2 = 383
$ws.Range("F23").Value2 = 130
$ws.Range("F24").Value2 = 3841
$ws.Range("F28").Value2 = 3042
$ws.Range("F29").Value2 = 2582
$ws.Range("F30").Value2 = 4404
$ws.Range("F32").Value2 = 942
$ws.Range("F33").Value2 = 2950
$ws.Range("F34").Value2 = 41
$ws.Range("F35").Value2 = 1064
$ws.Range("F36").Value2 = 1310
$ws.Range("F38").Value2 = 1160
$ws.Range("F39").Value2 = 750
$ws.Range("F41").Value2 = 453
$ws.Range("F44").Value2 = 32
$ws.Range("F45").Value2 = 23
$ws.Range("F46").Value2 = 153
$ws.Range("F47").Value2 = 333
$ws.Range("F48").Value2 = 3637
